$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 <- old row 8
$ws.Range('A5').Value = 'Stephen A. Klotz, Julia B. Jernberg, Richard Robbins'
$ws.Range('B5').Value = 'Department of Medicine, University of Arizona, Tucson; Department of Medicine, University of Arizona, Tucson; Department of Medicine, University of Arizona, Tucson'
$ws.Range('C5').Value = 'https://openalex.org/W4377826892'
$ws.Range('D5').Value = 'Turn Healthcare Workers Loose with Outpatient Telemedicine—Let Them Decide Its Fate; No Top-Down Decisions on What It Can and Cannot Do'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '2023-10-01'
$ws.Range('F5').Value = 'The American Journal of Medicine'
$ws.Range('G5').Value = 'Elsevier BV'
$ws.Range('H5').Value = 'https://doi.org/10.1016/j.amjmed.2023.05.005'
$ws.Range('I5').Value = 'N/A'
$ws.Range('J5').Value = 'publishedVersion'
$ws.Range('K5').Value = 'bronze'
$ws.Range('L5').Value = 'en'
$ws.Range('M5').Value = '0'
$ws.Range('N5').Value = '2023'
$ws.Range('O5').Value = 'https://pubmed.ncbi.nlm.nih.gov/37230402'
$ws.Range('P5').Value = 'https://doi.org/10.1016/j.amjmed.2023.05.005'
$ws.Range('Q5').Value = 'article'

# Row 6 <- old row 5
$ws.Range('A6').Value = 'Greta J. Binford, Samuel D. Robinson, Stephen A. Klotz'
$ws.Range('B6').Value = 'Professor of Biology, Lewis & Clark College, Portland, OR, 97219, USA; Institute for Molecular Bioscience, The University of Queensland, Brisbane, QLD, 4072, Australia; Division of Infectious Diseases, Department of Medicine, University of Arizona, 1501 N. Campbell Ave., Tucson, AZ, 85724, USA'
$ws.Range('C6').Value = 'https://openalex.org/W4386913058'
$ws.Range('D6').Value = 'Justin O Schmidt - His extraordinary impact on toxinology and arthropod biodiversity science'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '2023-10-01'
$ws.Range('F6').Value = 'Toxicon'
$ws.Range('G6').Value = 'Elsevier BV'
$ws.Range('H6').Value = 'https://doi.org/10.1016/j.toxicon.2023.107287'
$ws.Range('I6').Value = 'N/A'
$ws.Range('J6').Value = 'publishedVersion'
$ws.Range('K6').Value = 'bronze'
$ws.Range('L6').Value = 'en'
$ws.Range('M6').Value = '0'
$ws.Range('N6').Value = '2023'
$ws.Range('O6').Value = 'https://pubmed.ncbi.nlm.nih.gov/37740990'
$ws.Range('P6').Value = 'https://doi.org/10.1016/j.toxicon.2023.107287'
$ws.Range('Q6').Value = 'article'

# Row 7 <- old row 6
$ws.Range('A7').Value = 'Stephen A. Klotz, Krystal Fimbres, Lawrence York'
$ws.Range('B7').Value = '; ; '
$ws.Range('C7').Value = 'https://openalex.org/W4386987707'
$ws.Range('D7').Value = 'Infectious Diseases Telemedicine to the Arizona Department of Corrections During SARS-CoV-2 Pandemic. A Short Report.'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '2023-09-23'
$ws.Range('F7').Value = 'Southwest journal of pulmonary, critical care & sleep'
$ws.Range('G7').Value = 'N/A'
$ws.Range('H7').Value = 'https://doi.org/10.13175/swjpccs038-23'
$ws.Range('I7').Value = 'N/A'
$ws.Range('J7').Value = 'publishedVersion'
$ws.Range('K7').Value = 'bronze'
$ws.Range('L7').Value = 'en'
$ws.Range('M7').Value = '0'
$ws.Range('N7').Value = '2023'
$ws.Range('O7').Value = 'NA'
$ws.Range('P7').Value = 'https://doi.org/10.13175/swjpccs038-23'
$ws.Range('Q7').Value = 'article'

# Row 8 <- old row 7
$ws.Range('A8').Value = 'Nathan Kummet, Neha Mishra, Ana Claudia Marques Barbosa Diaz, Nicholas Cusick, Stephen A. Klotz, Nafees Ahmad'
$ws.Range('B8').Value = '; ; ; ; ; '
$ws.Range('C8').Value = 'https://openalex.org/W4387304862'
$ws.Range('D8').Value = 'Genetic Characterization of HIV-1 tat Gene from Virologically Controlled HIV-infected Older Patients on Long-term Antiretroviral Therapy'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '2023-10-02'
$ws.Range('F8').Value = 'N/A'
$ws.Range('G8').Value = 'N/A'
$ws.Range('H8').Value = 'https://doi.org/10.20944/preprints202310.0046.v1'
$ws.Range('I8').Value = 'N/A'
$ws.Range('J8').Value = 'submittedVersion'
$ws.Range('K8').Value = 'bronze'
$ws.Range('L8').Value = 'en'
$ws.Range('M8').Value = '0'
$ws.Range('N8').Value = '2023'
$ws.Range('O8').Value = 'NA'
$ws.Range('P8').Value = 'https://doi.org/10.20944/preprints202310.0046.v1'
$ws.Range('Q8').Value = 'article'
